$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.175812333333333
$ws.Range("H2").Value = 3.5274369999999999
$ws.Range("I2").Value = 0.033760122822238177
$ws.Range("J2").Value = 0.033760122822238177
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.66666666666666663
$ws.Range("M2").Value = 0.01848533333333334
$ws.Range("N2").Value = 0.055456000000000012
$ws.Range("O2").Value = 0.001625201930372746
$ws.Range("P2").Value = 0.001625201930372746
$ws.Range("Q2").Value = 0.021735282919111108
$ws.Range("R2").Value = 0.19561754627200001
$ws.Range("S2").Value = 0.00005486701678032249
$ws.Range("T2").Value = 0.00005486701678032248
$ws.Range("G3").Value = 1.175812333333333
$ws.Range("H3").Value = 3.5274369999999999
$ws.Range("I3").Value = 0.033760122822238177
$ws.Range("J3").Value = 0.033760122822238177
$ws.Range("O3").Value = 0.0026983345812381021
$ws.Range("P3").Value = 0.0026983345812381021
$ws.Range("Q3").Value = 0.036087248259777778
$ws.Range("R3").Value = 0.324785234338
$ws.Range("S3").Value = 0.00009109610687809097
$ws.Range("T3").Value = 0.00009109610687809094
$ws.Range("G4").Value = 1.175812333333333
$ws.Range("H4").Value = 3.5274369999999999
$ws.Range("I4").Value = 0.033760122822238177
$ws.Range("J4").Value = 0.033760122822238177
$ws.Range("M4").Value = 11.32499966666667
$ws.Range("N4").Value = 33.974998999999997
$ws.Range("O4").Value = 0.99567646348838923
$ws.Range("P4").Value = 0.99567646348838901
$ws.Range("Q4").Value = 13.31607428306256
$ws.Range("R4").Value = 119.844668547563
$ws.Range("S4").Value = 0.033614159698579768
$ws.Range("T4").Value = 0.033614159698579768
$ws.Range("I5").Value = 0.1056847291063769
$ws.Range("J5").Value = 0.1056847291063769
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.66666666666666663
$ws.Range("M5").Value = 0.01848533333333334
$ws.Range("N5").Value = 0.055456000000000012
$ws.Range("O5").Value = 0.001625201930372746
$ws.Range("P5").Value = 0.001625201930372746
$ws.Range("Q5").Value = 0.068041443434666685
$ws.Range("R5").Value = 0.61237299091200015
$ws.Range("S5").Value = 0.00017175902575460451
$ws.Range("T5").Value = 0.0001717590257546044
$ws.Range("I6").Value = 0.1056847291063769
$ws.Range("J6").Value = 0.1056847291063769
$ws.Range("O6").Value = 0.0026983345812381021
$ws.Range("P6").Value = 0.0026983345812381021
$ws.Range("S6").Value = 0.00028517275925651778
$ws.Range("T6").Value = 0.00028517275925651772
$ws.Range("I7").Value = 0.1056847291063769
$ws.Range("J7").Value = 0.1056847291063769
$ws.Range("M7").Value = 11.32499966666667
$ws.Range("N7").Value = 33.974998999999997
$ws.Range("O7").Value = 0.99567646348838923
$ws.Range("P7").Value = 0.99567646348838901
$ws.Range("Q7").Value = 41.685443823055337
$ws.Range("R7").Value = 375.16899440749808
$ws.Range("S7").Value = 0.10522779732136581
$ws.Range("T7").Value = 0.10522779732136581
$ws.Range("G8").Value = 29.971791333333329
$ws.Range("H8").Value = 89.915374
$ws.Range("I8").Value = 0.86055514807138489
$ws.Range("J8").Value = 0.86055514807138489
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.66666666666666663
$ws.Range("M8").Value = 0.01848533333333334
$ws.Range("N8").Value = 0.055456000000000012
$ws.Range("O8").Value = 0.001625201930372746
$ws.Range("P8").Value = 0.001625201930372746
$ws.Range("Q8").Value = 0.55403855339377783
$ws.Range("R8").Value = 4.9863469805440008
$ws.Range("S8").Value = 0.0013985758878378191
$ws.Range("T8").Value = 0.0013985758878378191
$ws.Range("G9").Value = 29.971791333333329
$ws.Range("H9").Value = 89.915374
$ws.Range("I9").Value = 0.86055514807138489
$ws.Range("J9").Value = 0.86055514807138489
$ws.Range("O9").Value = 0.0026983345812381021
$ws.Range("P9").Value = 0.0026983345812381021
$ws.Range("Q9").Value = 0.9198742384084444
$ws.Range("R9").Value = 8.2788681456760003
$ws.Range("S9").Value = 0.0023220657151034928
$ws.Range("T9").Value = 0.0023220657151034928
$ws.Range("G10").Value = 29.971791333333329
$ws.Range("H10").Value = 89.915374
$ws.Range("I10").Value = 0.86055514807138489
$ws.Range("J10").Value = 0.86055514807138489
$ws.Range("M10").Value = 11.32499966666667
$ws.Range("N10").Value = 33.974998999999997
$ws.Range("O10").Value = 0.99567646348838923
$ws.Range("P10").Value = 0.99567646348838901
$ws.Range("Q10").Value = 339.43052685940302
$ws.Range("R10").Value = 3054.8747417346258
$ws.Range("S10").Value = 0.85683450646844361
$ws.Range("T10").Value = 0.8568345064684435
